$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: bulletin Volume/Number and the reporting week's date range.
# These live inside rich-text shared strings referenced by A8 and C8.
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(16, 2).Text = "28"
$ws.Range("C8").Characters(27, 9).Text = "7/7/2025"
$ws.Range("C8").Characters(48, 8).Text = "7/13/2025"

# ---------------------------------------------------------------------------
# Row 14
# ---------------------------------------------------------------------------
$ws.Range("N14").Value = -91.666666666666

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
$ws.Range("D15").Value = 2
$ws.Range("C15").Copy($ws.Range("F15")) | Out-Null
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = -46.153846153846
$ws.Range("N15").Value = -70.833333333333

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("I16").Value = 45
$ws.Range("J16").Value = 27
$ws.Range("K16").Value = 66.666666666666
$ws.Range("L16").Value = -13.461538461538
$ws.Range("M16").Value = -74.285714285714
$ws.Range("N16").Value = -91.362763915547

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -47.619047619047
$ws.Range("I17").Value = 110
$ws.Range("J17").Value = 147
$ws.Range("K17").Value = -25.170068027210
$ws.Range("L17").Value = -13.385826771653
$ws.Range("M17").Value = -32.098765432098
$ws.Range("N17").Value = -44.162436548223

# ---------------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 180
$ws.Range("I18").Value = 76
$ws.Range("J18").Value = 58
$ws.Range("K18").Value = 31.034482758620
$ws.Range("L18").Value = 2.702702702702
$ws.Range("M18").Value = -58.918918918918
$ws.Range("N18").Value = -91.090269636576

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 6
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 25
$ws.Range("I19").Value = 191
$ws.Range("J19").Value = 198
$ws.Range("K19").Value = -3.535353535353
$ws.Range("L19").Value = -8.173076923076
$ws.Range("M19").Value = -18.376068376068
$ws.Range("N19").Value = -39.747634069400

# ---------------------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------------------
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 60
$ws.Range("I20").Value = 100
$ws.Range("J20").Value = 122
$ws.Range("K20").Value = -18.032786885245
$ws.Range("L20").Value = 19.047619047619
$ws.Range("M20").Value = -48.979591836734
$ws.Range("N20").Value = -94.378864530635

# ---------------------------------------------------------------------------
# Row 21
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -5
$ws.Range("F21").Value = 71
$ws.Range("G21").Value = 70
$ws.Range("H21").Value = 1.428571428571
$ws.Range("I21").Value = 530
$ws.Range("J21").Value = 565
$ws.Range("K21").Value = -6.194690265486
$ws.Range("L21").Value = -4.159132007233
$ws.Range("M21").Value = -45.473251028806
$ws.Range("N21").Value = -85.687280583310

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 12
$ws.Range("E24").Value = -14.285714285714
$ws.Range("F24").Value = 59
$ws.Range("G24").Value = 41
$ws.Range("H24").Value = 43.902439024390
$ws.Range("I24").Value = 311
$ws.Range("J24").Value = 375
$ws.Range("K24").Value = -17.066666666666
$ws.Range("L24").Value = -25.598086124401
$ws.Range("M24").Value = -30.734966592427

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = -18.75
$ws.Range("I25").Value = 80
$ws.Range("J25").Value = 95
$ws.Range("K25").Value = -15.789473684210
$ws.Range("L25").Value = 2.564102564102

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 28
$ws.Range("G26").Value = 54
$ws.Range("H26").Value = -48.148148148148
$ws.Range("I26").Value = 233
$ws.Range("J26").Value = 229
$ws.Range("K26").Value = 1.746724890829
$ws.Range("L26").Value = 28.729281767955
$ws.Range("M26").Value = -34.733893557423

# ---------------------------------------------------------------------------
# Row 27
# ---------------------------------------------------------------------------
$ws.Range("D27").Value = 2
$ws.Range("C15").Copy($ws.Range("F27")) | Out-Null
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -61.111111111111

# ---------------------------------------------------------------------------
# Row 28 (D28/E28 flip from the "n/a" placeholder text back to real numbers)
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = 1
$ws.Range("F27").Copy($ws.Range("D28")) | Out-Null
$ws.Range("D28").Value = 1
$ws.Range("H28").Copy($ws.Range("E28")) | Out-Null
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 20
$ws.Range("J28").Value = 24
$ws.Range("K28").Value = -16.666666666666
$ws.Range("L28").Value = 66.666666666666

# ---------------------------------------------------------------------------
# Row 29 (D29/E29 flip from real numbers to the "n/a" placeholder text)
# ---------------------------------------------------------------------------
$ws.Range("C22").Copy($ws.Range("D29")) | Out-Null
$ws.Range("E22").Copy($ws.Range("E29")) | Out-Null

# ---------------------------------------------------------------------------
# Row 30 (same transformation as row 29)
# ---------------------------------------------------------------------------
$ws.Range("C22").Copy($ws.Range("D30")) | Out-Null
$ws.Range("E22").Copy($ws.Range("E30")) | Out-Null

# ---------------------------------------------------------------------------
# Row 33
# ---------------------------------------------------------------------------
$ws.Range("C22").Copy($ws.Range("C33")) | Out-Null
